$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Value = 0.90076854960269
$ws.Range("F1").Value = -1.570796383038867

$ws.Range("E2").Value = 0.8988340247709597
$ws.Range("F2").Value = -1.570796377815036

$ws.Range("E3").Value = 0.8901660709669087
$ws.Range("F3").Value = -1.570796354408809

$ws.Range("E4").Value = 0.877971082939261
$ws.Range("F4").Value = -1.570796321478471

$ws.Range("E5").Value = 0.8693031291352099
$ws.Range("F5").Value = -1.570796298072244

$ws.Range("E6").Value = 0.8673686043034796
$ws.Range("F6").Value = -1.570796292848413
